$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel alignment constants (used throughout)
$xlLeft    = -4131
$xlCenter  = -4108
$xlBottom  = -4107
$xlGeneral = 1
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Rows 4-6: the "ID" cells (B4,B5,B6) gain vertical centering (they already
#    had horizontal centering + the right fill/border).
# ---------------------------------------------------------------------------
$ws.Range("B4").VerticalAlignment = $xlCenter
$ws.Range("B5").VerticalAlignment = $xlCenter
$ws.Range("B6").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# 2. Row 7: description cell (C7) loses its special "shadow" font, its wrap,
#    and its left/vertical-center alignment -- it becomes a plain fill-3 cell
#    like B7 -- and its text is replaced by the new story.
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)
$ws.Range("C7").HorizontalAlignment = $xlGeneral
$ws.Range("C7").VerticalAlignment = $xlBottom
$ws.Range("C7").WrapText = $false
$ws.Range("C7").Value = "Eu, como funcionário de Infraestrutura, gostaria de ter acesso a gráficos para melho visulização"

# ---------------------------------------------------------------------------
# 3. Row 8 (new): same look as rows 4/6 for the ID column, same look as row 4
#    for the description column minus the vertical centering.
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Value = 5

$ws.Range("C4").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("C8").VerticalAlignment = $xlBottom
$ws.Range("C8").Value = "Eu, como desenvolvedor, gostaria de implementar DataMining para alcançar uma margem de erro mais exata"

# ---------------------------------------------------------------------------
# 4. Row 9 (new): ID cell looks like row 5, description is a fill-3 cell with
#    only wrap turned on.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)
$ws.Range("B9").Value = 6

$ws.Range("C7").Copy()
$ws.Range("C9").PasteSpecial($xlPasteFormats)
$ws.Range("C9").WrapText = $true
$ws.Range("C9").Value = "Eu, como funcionário de Infraestrutura, quero ver um histórico de funcionamento das máquinas para uma melhor administração dos Totens"

# ---------------------------------------------------------------------------
# 5. Row 10 (new): ID cell looks like row 4/6/8, description is a fill-4 cell
#    with only wrap turned on (no horizontal/vertical alignment).
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B10").PasteSpecial($xlPasteFormats)
$ws.Range("B10").Value = 7

$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial($xlPasteFormats)
$ws.Range("C10").HorizontalAlignment = $xlGeneral
$ws.Range("C10").VerticalAlignment = $xlBottom
$ws.Range("C10").Value = "Eu, como gerente, gostaria de utilizar o sistema em mais de um totem para uma melhor gestão de tempo e custos"

# ---------------------------------------------------------------------------
# 6. Row 11 (new): ID cell looks like row 7, description is a fill-3 cell with
#    only wrap turned on (same shape as row 9's description). Its shared
#    string is registered AFTER row 13's (matches the author's entry order).
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("B11").Value = 8

$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 7. Row 12 (new): ID cell looks like row 4/6/8/10, description reuses the
#    "shadow font" look from row 6 but on fill 4 without the wrap, and holds
#    the old "site institucional" story that used to live in row 7.
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B12").PasteSpecial($xlPasteFormats)
$ws.Range("B12").Value = 9

$ws.Range("C6").Copy()
$ws.Range("C12").PasteSpecial($xlPasteFormats)
$ws.Range("C12").WrapText = $false
$ws.Range("C12").Value = "Eu como desenvolvedor, gostaria de um site Institucional para atigir um maior público. "

# ---------------------------------------------------------------------------
# 8. Row 13 (new): ID cell looks like row 5/9, description is a fill-3 cell
#    with only wrap turned on (same shape as rows 9/11's description). Its
#    text is registered in the shared-string table before row 11's text.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)
$ws.Range("B13").Value = 10

$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Value = "Eu, como desenvolvedor, gostaria de fazer uma solução em Java para ter mais conforto em relação a tecnologia"

# Now that "Java" has claimed its shared-string slot, give row 11 its text.
$ws.Range("C11").Value = "Eu, como fucionário de Infraestrutura, preciso receber Logs para poder me informar sobre os incidentes"

# ---------------------------------------------------------------------------
# 9. Row heights: new wrapped rows render two lines tall like the existing
#    wrapped rows (4,5,6); the non-wrapped rows (7,12) stay single-line.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30

# ---------------------------------------------------------------------------
# 10. Selection moves to E8 (matches the saved sheetView in the target file).
# ---------------------------------------------------------------------------
$ws.Range("E8").Select()
